$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Ideas"

$ws.Range("A1").Value = "Thematic"
$ws.Range("C1").Value = "Function"
$ws.Range("E1").Value = "Mechanics"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true

Write-Host "done"
